$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 3 ("Before legends is a round based strategy-RPG fusion ..."):
# split into a "Vision 01" paragraph followed by a new paragraph holding the
# original sentence.
# ---------------------------------------------------------------------------
$pOrig = $d.Paragraphs.Item(3)
$origText = $pOrig.Range.Text
$origText = $origText.TrimEnd([char]13, [char]7)

$pOrig.Range.InsertParagraphAfter()

$pVision01 = $d.Paragraphs.Item(3)
$pVision01.Range.Text = "Vision 01"
$pVision01.Range.LanguageID = "en-US"

$pBefore = $d.Paragraphs.Item(4)
$pBefore.Range.Text = $origText
$pBefore.Range.LanguageID = "en-US"

# ---------------------------------------------------------------------------
# Paragraph 5 (originally empty, right before "key words:"): becomes
# "Vision 02", followed by the new mission paragraph, "Vision 03" (carrying
# the relocated _GoBack bookmark) and a new empty paragraph.
# ---------------------------------------------------------------------------
$pEmpty = $d.Paragraphs.Item(5)
$pEmpty.Range.Text = "Vision 02"
$pEmpty.Range.LanguageID = "en-US"

# New paragraph with the revised vision text (three runs, matching the source).
$pEmpty.Range.InsertParagraphAfter()
$pVision02Body = $d.Paragraphs.Item(6)
$run1 = "Before Legends is a round based strategy PC-game with RPG elements, set in a 3D prehistorical fantasy environment, in which the player will take charge of a small tribe and lead it against the challenges of nature and other tribes "
$run2 = "competing"
$run3 = " for survival on the journey to forge their first legend and become a civilization."
$pVision02Body.Range.Text = $run1
$pVision02Body.Range.LanguageID = "en-US"
$pVision02Body.Range.InsertAfter($run2)
$pVision02Body.Range.InsertAfter($run3)

# New paragraph "Vision 03" (will carry the relocated bookmark).
$pVision02Body.Range.InsertParagraphAfter()
$pVision03 = $d.Paragraphs.Item(7)
$pVision03.Range.Text = "Vision 03X"
$pVision03.Range.LanguageID = "en-US"

# Place the _GoBack bookmark right after "Vision 03" (before the trailing
# placeholder "X"); placing a zero-width range exactly on a paragraph-mark
# character is unreliable, so we insert a throwaway character, bookmark just
# before it, then delete the throwaway character.
$rSafe = $d.Paragraphs.Item(7).Range
$safePos = $rSafe.End - 2
$safeRange = $d.Range($safePos, $safePos)

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$d.Bookmarks.Add("_GoBack", $safeRange)

$rTrim = $d.Paragraphs.Item(7).Range
$trimPos = $rTrim.End - 1
$trimRange = $d.Range($trimPos - 1, $trimPos)
$trimRange.Text = ""

# New empty paragraph after "Vision 03".
$d.Paragraphs.Item(7).Range.InsertParagraphAfter()
